$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.169.49'
$ws.Range("E2").Value = '  -3.07%  '

# Row 3
$ws.Range("D3").Value = '3.206.10'
$ws.Range("E3").Value = '  -2.34%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.47'
$ws.Range("E5").Value = '  -2.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.22'
$ws.Range("E6").Value = '  -5.54%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  -4.53%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").Value = '3.209.66'
$ws.Range("E9").Value = '  -2.18%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -2.78%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.79'
$ws.Range("E11").Value = '  +0.67%  '

# Row 12
$ws.Range("E12").Value = '  -4.20%  '

# Row 13
$ws.Range("D13").Value = '3.768.56'
$ws.Range("E13").Value = '  -2.27%  '

# Row 14
$ws.Range("E14").Value = '  -1.35%  '

# Row 15
$ws.Range("D15").Value = '64.277.09'
$ws.Range("E15").Value = '  -2.92%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.63'
$ws.Range("E16").Value = '  -3.25%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000160'
$ws.Range("E17").Value = '  -2.33%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.212.92'
$ws.Range("E18").Value = '  -1.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '418.07'
$ws.Range("E19").Value = '  -3.94%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.41'
$ws.Range("E20").Value = '  -1.80%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.89'
$ws.Range("E21").Value = '  -2.65%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.23'
$ws.Range("E22").Value = '  -2.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.22'
$ws.Range("E24").Value = '  -1.23%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.67'
$ws.Range("E25").Value = '  -1.13%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.204'
$ws.Range("E26").Value = '  +3.99%  '

# Row 27
$ws.Range("E27").Value = '  -3.08%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000111'
$ws.Range("E28").Value = '  -2.55%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.89'
$ws.Range("E29").Value = '  +0.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.08%  '

# Row 31
$ws.Range("E31").Value = '  -5.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.99'
$ws.Range("E32").Value = '  -1.48%  '

# Row 33
$ws.Range("E33").Value = '  +0.08%  '

# Row 34
$ws.Range("E34").Value = '  -3.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.45'
$ws.Range("E35").Value = '  -2.94%  '

# Row 36
$ws.Range("E36").Value = '  -4.10%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.11'
$ws.Range("E37").Value = '  -0.63%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.40'
$ws.Range("E38").Value = '  -1.69%  '

# Row 39
$ws.Range("D39").Value = '2.763.56'
$ws.Range("E39").Value = '  -1.11%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.72'
$ws.Range("E40").Value = '  -3.29%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.40'
$ws.Range("E41").Value = '  -3.53%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.23'
$ws.Range("E42").Value = '  -2.87%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.721'
$ws.Range("E43").Value = '  -6.69%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.21'
$ws.Range("E44").Value = '  -2.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.81'
$ws.Range("E45").Value = '  -5.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0630'
$ws.Range("E46").Value = '  -4.80%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.19'
$ws.Range("E47").Value = '  -4.77%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '300.50'
$ws.Range("E48").Value = '  -6.00%  '

# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.15'
$ws.Range("E49").Value = '  -7.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0263'
$ws.Range("E50").Value = '  -1.96%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("E51").Value = '  -4.33%  '
